$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update column C ("Förändrad") from row 2 to row 27 with the new date serial value (45324 = 2024-02-02)
for ($r = 2; $r -le 27; $r++) {
    $ws.Cells.Item($r, 3).Value = 45324
}
